# Update crypto price/volume table per latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to remain plain text (matches the workbook's existing
    # inline-string cells) even when the new value looks numeric, e.g. "587.66"
    # or a pseudo price like "64.932.32". Restore the default "Normal" style
    # afterwards so we do not leave a stray text-number-format on the cell.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "64.932.32"
Set-TextCell "E2" "  -0.47%  "
Set-TextCell "D3" "3.524.57"
Set-TextCell "E3" "  -1.24%  "
Set-TextCell "E4" "  -0.10%  "
Set-TextCell "D5" "587.66"
Set-TextCell "E5" "  -2.22%  "
Set-TextCell "D6" "134.02"
Set-TextCell "E6" "  -1.54%  "
Set-TextCell "D7" "3.524.05"
Set-TextCell "E7" "  -1.22%  "
Set-TextCell "D9" "0.491"
Set-TextCell "E9" "  -1.05%  "
Set-TextCell "D10" "0.125"
Set-TextCell "E10" "  +1.07%  "
Set-TextCell "D11" "7.17"
Set-TextCell "E11" "  +2.20%  "
Set-TextCell "D12" "0.386"
Set-TextCell "E12" "  -0.55%  "
Set-TextCell "D13" "4.123.89"
Set-TextCell "E13" "  -1.26%  "
Set-TextCell "D14" "27.78"
Set-TextCell "E14" "  +1.99%  "
Set-TextCell "D15" "0.0000181"
Set-TextCell "E15" "  -1.41%  "
Set-TextCell "E16" "  +0.61%  "
Set-TextCell "D17" "3.524.18"
Set-TextCell "E17" "  -1.55%  "
Set-TextCell "D18" "64.932.45"
Set-TextCell "E18" "  +0.36%  "
Set-TextCell "D19" "10.07"
Set-TextCell "E19" "  +0.21%  "
Set-TextCell "D20" "14.26"
Set-TextCell "E20" "  -1.06%  "
Set-TextCell "D21" "5.68"
Set-TextCell "E21" "  -3.24%  "
Set-TextCell "D22" "390.92"
Set-TextCell "E22" "  -0.85%  "
Set-TextCell "D23" "0.578"
Set-TextCell "E23" "  -0.47%  "
Set-TextCell "D24" "3.664.00"
Set-TextCell "E24" "  -1.43%  "
Set-TextCell "D25" "74.48"
Set-TextCell "E25" "  +0.58%  "
Set-TextCell "E26" "  +0.21%  "
Set-TextCell "D27" "0.0000111"
Set-TextCell "E27" "  -2.78%  "
Set-TextCell "D28" "1.62"
Set-TextCell "E28" "  +4.54%  "
Set-TextCell "B29" "RenderToken"
Set-TextCell "C29" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D29" "7.57"
Set-TextCell "E29" "  -2.81%  "
Set-TextCell "B30" "Binance-PegBSC-USD"
Set-TextCell "C30" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D30" "0.995"
Set-TextCell "E30" "  -0.31%  "
Set-TextCell "D31" "2.27"
Set-TextCell "E31" "  -1.56%  "
Set-TextCell "D32" "8.31"
Set-TextCell "E32" "  -2.85%  "
Set-TextCell "D33" "3.526.57"
Set-TextCell "B34" "EthereumClassic"
Set-TextCell "C34" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D34" "24.08"
Set-TextCell "E34" "  -0.40%  "
Set-TextCell "B35" "USDe"
Set-TextCell "C35" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D35" "1.00"
Set-TextCell "E35" "  +0.04%  "
Set-TextCell "D36" "0.146"
Set-TextCell "E36" "  +0.12%  "
Set-TextCell "B37" "ImmutableX"
Set-TextCell "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D37" "1.60"
Set-TextCell "E37" "  +1.30%  "
Set-TextCell "D38" "5.24"
Set-TextCell "E38" "  +3.43%  "
Set-TextCell "B39" "Monero"
Set-TextCell "C39" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D39" "172.12"
Set-TextCell "E39" "  +1.83%  "
Set-TextCell "B40" "Aptos"
Set-TextCell "C40" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D40" "6.98"
Set-TextCell "E40" "  -0.46%  "
Set-TextCell "D41" "0.0813"
Set-TextCell "E41" "  +0.16%  "
Set-TextCell "D42" "0.820"
Set-TextCell "E42" "  -0.94%  "
Set-TextCell "D43" "26.45"
Set-TextCell "E43" "  -0.17%  "
Set-TextCell "B44" "ONDO"
Set-TextCell "C44" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell "D44" "1.26"
Set-TextCell "E44" "  +1.32%  "
Set-TextCell "B45" "FirstDigitalUSD"
Set-TextCell "C45" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D45" "1.00"
Set-TextCell "E45" "  -0.11%  "
Set-TextCell "B46" "OKB"
Set-TextCell "C46" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D46" "42.45"
Set-TextCell "E46" "  -1.45%  "
Set-TextCell "D47" "4.43"
Set-TextCell "E47" "  -0.99%  "
Set-TextCell "D48" "1.67"
Set-TextCell "E48" "  -0.33%  "
Set-TextCell "D49" "2.474.73"
Set-TextCell "E49" "  -0.22%  "
Set-TextCell "D50" "6.88"
Set-TextCell "E50" "  -0.62%  "
Set-TextCell "D51" "0.900"
Set-TextCell "E51" "  +2.04%  "
